$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Whitelist now only keeps two Steam64 ids (rows 2 & 3); the old sample
# 1..12 placeholder list is replaced by the real allowed ids.
$ws.Range("A2").Value = 21942357
$ws.Range("A3").Value = 170633010

# Drop the old rows 4:16 (data rows 4:13 + the trailing blank rows 14:16)
# so the sheet's used range shrinks back down to A2:A3.
$ws.Range("A4:A16").EntireRow.Delete()

# Leave the cursor parked one row below the data, matching where the
# author left the selection after entering the new ids.
$ws.Range("B4").Select() | Out-Null
